$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.487.87'
$ws.Range("E2").Value = '  -1.35%  '

$ws.Range("D3").Value = '2.356.47'
$ws.Range("E3").Value = '  -0.71%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = '321.04'
$ws.Range("E5").Value = '  -0.68%  '

$ws.Range("D6").Value = '107.27'
$ws.Range("E6").Value = '  +3.52%  '

$ws.Range("E7").Value = '  -1.09%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  -5.61%  '

$ws.Range("D10").Value = '41.30'
$ws.Range("E10").Value = '  -0.14%  '

$ws.Range("D11").Value = '0.0923'
$ws.Range("E11").Value = '  -1.31%  '

$ws.Range("D12").Value = '8.47'
$ws.Range("E12").Value = '  -0.71%  '

$ws.Range("D13").Value = '0.995'
$ws.Range("E13").Value = '  -2.83%  '

$ws.Range("E14").Value = '  +0.47%  '

$ws.Range("D15").Value = '15.98'
$ws.Range("E15").Value = '  -5.18%  '

$ws.Range("D16").Value = '2.713.33'
$ws.Range("E16").Value = '  -0.98%  '

$ws.Range("D17").Value = '2.340.72'
$ws.Range("E17").Value = '  -1.62%  '

$ws.Range("D18").Value = '42.510.64'
$ws.Range("E18").Value = '  -1.38%  '

$ws.Range("D19").Value = '7.64'
$ws.Range("E19").Value = '  -2.67%  '

$ws.Range("E20").Value = '  -1.38%  '

$ws.Range("D21").Value = '76.30'
$ws.Range("E21").Value = '  -0.19%  '

$ws.Range("D22").Value = '3.61'
$ws.Range("E22").Value = '  +6.98%  '

$ws.Range("D23").Value = '256.94'
$ws.Range("E23").Value = '  -7.09%  '

$ws.Range("D24").Value = '2.31'
$ws.Range("E24").Value = '  -3.27%  '

$ws.Range("D25").Value = '9.37'
$ws.Range("E25").Value = '  -1.36%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("D27").Value = '11.41'
$ws.Range("E27").Value = '  -1.93%  '

$ws.Range("D28").Value = '22.81'
$ws.Range("E28").Value = '  -1.39%  '

$ws.Range("E29").Value = '  +3.02%  '

$ws.Range("D30").Value = '175.84'
$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("D31").Value = '36.68'
$ws.Range("E31").Value = '  -3.09%  '

$ws.Range("D32").Value = '0.0889'
$ws.Range("E32").Value = '  -2.88%  '

$ws.Range("D33").Value = '6.05'
$ws.Range("E33").Value = '  +3.70%  '

$ws.Range("D34").Value = '2.88'
$ws.Range("E34").Value = '  -9.71%  '

$ws.Range("D35").Value = '0.127'
$ws.Range("E35").Value = '  +19.67%  '

$ws.Range("E36").Value = '  -1.06%  '

$ws.Range("D37").Value = '4.62'
$ws.Range("E37").Value = '  -4.64%  '

$ws.Range("D38").Value = '0.0363'
$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("D39").Value = '3.85'
$ws.Range("E39").Value = '  -7.85%  '

$ws.Range("D40").Value = '2.68'
$ws.Range("E40").Value = '  -4.35%  '

$ws.Range("D41").Value = '0.239'
$ws.Range("E41").Value = '  +4.02%  '

$ws.Range("D42").Value = '71.19'
$ws.Range("E42").Value = '  +2.95%  '

$ws.Range("D43").Value = '1.47'
$ws.Range("E43").Value = '  -6.77%  '

$ws.Range("E44").Value = '  -0.27%  '

$ws.Range("D45").Value = '12.01'
$ws.Range("E45").Value = '  -2.78%  '

$ws.Range("D46").Value = '112.65'
$ws.Range("E46").Value = '  -9.06%  '

$ws.Range("D47").Value = '5.48'
$ws.Range("E47").Value = '  -1.43%  '

$ws.Range("D48").Value = '9.11'
$ws.Range("E48").Value = '  -3.86%  '

$ws.Range("D49").Value = '84.76'
$ws.Range("E49").Value = '  -8.57%  '

$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").Value = '1.28'
$ws.Range("E50").Value = '  -1.11%  '

$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").Value = '73.93'
$ws.Range("E51").Value = '  +3.06%  '
